# Actualización automática 2025-06-19 16:40:08
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M5").Value = 2758.82
$ws1.Range("C19").Value = 518.4
$ws1.Range("E29").Value = 64.81999999999999
$ws1.Range("C53").Value = "6 de 51"
$ws1.Range("E53").Value = "2 de 51"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F5").Value = 2760.94
$ws2.Range("F19").Value = 518.4
$ws2.Range("F29").Value = 92.47
$ws2.Range("F53").Value = 43508.04

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D2").Value = 3576.95
$ws3.Range("E2").Value = 6393.39304517915
$ws3.Range("F2").Value = 0.3587589698560595

$ws3.Range("D4").Value = 641.34
$ws3.Range("E4").Value = 361.66
$ws3.Range("F4").Value = 0.6394217347956132

$ws3.Range("D16").Value = 9659.459999999999
$ws3.Range("E16").Value = 23081.99
$ws3.Range("F16").Value = 0.2950223646173276

$ws3.Range("D19").Value = 43508.04
$ws3.Range("E19").Value = 50939.40064517915
$ws3.Range("F19").Value = 0.4606587505473159
